$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 'Estados Unidos'
$ws.Range("B4").Value = 2103189
$ws.Range("C4").Value = 13488
$ws.Range("D4").Value = 819375
$ws.Range("E4").Value = 1167377
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 403
$ws.Range("H4").Value = 116437

$ws.Range("A12").Value = 'Alemania'
$ws.Range("B12").Value = 187009
$ws.Range("C12").Value = 214
$ws.Range("D12").Value = 171600
$ws.Range("E12").Value = 6554
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 4
$ws.Range("H12").Value = 8855

$ws.Range("A16").Value = 'Francia'
$ws.Range("B16").Value = 156287
$ws.Range("C16").Value = 726
$ws.Range("D16").Value = 72149
$ws.Range("E16").Value = 54764
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 28
$ws.Range("H16").Value = 29374

$ws.Range("A96").Value = 'Mayotte'
$ws.Range("B96").Value = 2268
$ws.Range("C96").Value = 28
$ws.Range("D96").Value = 1790
$ws.Range("E96").Value = 450
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 28

$ws.Range("A97").Value = 'Croacia'
$ws.Range("B97").Value = 2249
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 2133
$ws.Range("E97").Value = 9
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 107

$ws.Range("A103").Value = 'Sri Lanka'
$ws.Range("B103").Value = 1880
$ws.Range("C103").Value = 3
$ws.Range("D103").Value = 1196
$ws.Range("E103").Value = 673
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 11

$ws.Range("A117").Value = 'Zambia'
$ws.Range("B117").Value = 1321
$ws.Range("C117").Value = 121
$ws.Range("D117").Value = 1104
$ws.Range("E117").Value = 207
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 10

$ws.Range("A118").Value = 'Guinea Ecuatorial'
$ws.Range("B118").Value = 1306
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 200
$ws.Range("E118").Value = 1094
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 12

$ws.Range("A119").Value = 'Paraguay'
$ws.Range("B119").Value = 1254
$ws.Range("C119").Value = 24
$ws.Range("D119").Value = 633
$ws.Range("E119").Value = 610
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 11

$ws.Range("A120").Value = 'Madagascar'
$ws.Range("B120").Value = 1240
$ws.Range("C120").Value = 37
$ws.Range("D120").Value = 344
$ws.Range("E120").Value = 886
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 10

$ws.Range("A126").Value = 'Republica de Chipre'
$ws.Range("B126").Value = 980
$ws.Range("C126").Value = 5
$ws.Range("D126").Value = 807
$ws.Range("E126").Value = 155
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 18

$ws.Range("A131").Value = 'Republica del Chad'
$ws.Range("B131").Value = 848
$ws.Range("C131").Value = 0
$ws.Range("D131").Value = 711
$ws.Range("E131").Value = 65
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 72

$ws.Range("A144").Value = 'Ruanda'
$ws.Range("B144").Value = 510
$ws.Range("C144").Value = 16
$ws.Range("D144").Value = 321
$ws.Range("E144").Value = 187
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 2

$ws.Range("A145").Value = 'Mozambique'
$ws.Range("B145").Value = 509
$ws.Range("C145").Value = 20
$ws.Range("D145").Value = 145
$ws.Range("E145").Value = 362
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 2

$ws.Range("A146").Value = 'Tanzania'
$ws.Range("B146").Value = 509
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 183
$ws.Range("E146").Value = 305
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 21

$ws.Range("A152").Value = 'Liberia'
$ws.Range("B152").Value = 421
$ws.Range("C152").Value = 11
$ws.Range("D152").Value = 210
$ws.Range("E152").Value = 179
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 1
$ws.Range("H152").Value = 32

$ws.Range("A155").Value = 'Zimbabue'
$ws.Range("B155").Value = 343
$ws.Range("C155").Value = 11
$ws.Range("D155").Value = 51
$ws.Range("E155").Value = 288
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 4

$ws.Range("A156").Value = 'Mauricio'
$ws.Range("B156").Value = 337
$ws.Range("C156").Value = 0
$ws.Range("D156").Value = 325
$ws.Range("E156").Value = 2
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 10

$ws.Range("A157").Value = 'Isla de Man'
$ws.Range("B157").Value = 336
$ws.Range("C157").Value = 0
$ws.Range("D157").Value = 312
$ws.Range("E157").Value = 0
$ws.Range("F157").Value = 0
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 24

$ws.Range("A158").Value = 'Vietnam'
$ws.Range("B158").Value = 333
$ws.Range("C158").Value = 1
$ws.Range("D158").Value = 323
$ws.Range("E158").Value = 10
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 0

$ws.Range("A166").Value = 'Guadalupe'
$ws.Range("B166").Value = 171
$ws.Range("C166").Value = 7
$ws.Range("D166").Value = 157
$ws.Range("E166").Value = 0
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 14

$ws.Range("A167").Value = 'Surinam'
$ws.Range("B167").Value = 168
$ws.Range("C167").Value = 0
$ws.Range("D167").Value = 9
$ws.Range("E167").Value = 157
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 2

$ws.Range("A168").Value = 'Siria'
$ws.Range("B168").Value = 164
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 68
$ws.Range("E168").Value = 90
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 6
